$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-196 down to 71-197
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new data record
$ws.Cells.Item(70,1).Value = 10
$ws.Cells.Item(70,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70,3).Value = "La Araucanía"
$ws.Cells.Item(70,4).Value = 44469
$ws.Cells.Item(70,5).Value = 9
$ws.Cells.Item(70,6).Value = 100112044
$ws.Cells.Item(70,7).Value = "Perejil"
$ws.Cells.Item(70,8).Value = "Sin especificar"
$ws.Cells.Item(70,9).Value = "Primera"
$ws.Cells.Item(70,10).Value = 60
$ws.Cells.Item(70,11).Value = 4000
$ws.Cells.Item(70,12).Value = 5000
$ws.Cells.Item(70,13).Value = 4500
$ws.Cells.Item(70,14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(70,15).Value = "Provincia de Cautín"
$ws.Cells.Item(70,16).Value = 1500
$ws.Cells.Item(70,17).Value = 3
$ws.Cells.Item(70,18).Value = "Hortaliza"
